$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "delete_past_homework" / "delete_empty_rows": the assignment that used to
# live in row 4 (WebDev / assignment 1) has been cleared out, leaving a
# blank spacer row behind it and pushing the remaining homework rows down.
$ws.Rows("4:4").Insert()

# Restore the focus/selection Excel left the workbook in after the edit.
$ws.Range("D19").Select()
